$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 390 (shifts the existing 390-423 down to 392-425)
$ws.Rows.Item(390).Insert()
$ws.Rows.Item(390).Insert()

# New weekly records to insert (columns A-C repeat the market/region for every row)
$newRows = @(
    @(11, "Vega Monumental Concepción", "Bíobío", 45013, 8, 100114001, "Papa", "Asterix",   "1a (cosecha)",        250, 10000, 11000, 10400, "$/saco 25 kilos",  "Región de Los Lagos", 416, 25, "Hortaliza"),
    @(11, "Vega Monumental Concepción", "Bíobío", 45013, 8, 100114001, "Papa", "Patagonia", "1a (cosecha lavada)", 220, 10000, 11000, 10545, "$/malla 25 kilos", "Región de Los Lagos", 422, 25, "Hortaliza")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 390 + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
